$d = $word.ActiveDocument

# 1. Update the cached date field result text: 7/1/2025 -> 7/23/2025
$d.Content.Find.Execute("7/1/2025", $false, $false, $false, $false, $false, $true, 1, $false, "7/23/2025", 2)

# 2. Turn the space inside the {{Client Name}} placeholder into an underscore
$d.Content.Find.Execute("{{Client Name}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{Client_Name}}", 2)

# 3. Turn the space inside the {{Brief Synopsis}} placeholder into an underscore
$d.Content.Find.Execute("{{Brief Synopsis}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{Brief_Synopsis}}", 2)

# 4. Turn the space inside the {{Settlement Demand}} placeholder into an underscore
$d.Content.Find.Execute("{{Settlement Demand}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{Settlement_Demand}}", 2)

Write-Host "done"
